$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update attendee names and attendance status for the remaining rows
$ws.Range("B2").Value = "DavidLondoño"
$ws.Range("C2").Value = "Presente"

$ws.Range("B3").Value = "CarlosRiaño"
$ws.Range("C3").Value = "Presente"

$ws.Range("B4").Value = "AndrésPerea"

# Remove the now-unneeded trailing rows (5 through 11)
$ws.Range("A5:A11").EntireRow.Delete()
